$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Header block: "Phone number: ..." line had six tabs before "Email:";
#    drop one of them so five tabs remain.
# ---------------------------------------------------------------------------
$find1 = $d.Content.Find
$find1.ClearFormatting()
$old1 = "718-612-993" + [char]9 + [char]9 + [char]9 + [char]9 + [char]9 + [char]9 + "Email"
$new1 = "718-612-993" + [char]9 + [char]9 + [char]9 + [char]9 + [char]9 + "Email"
$find1.Execute($old1, $false, $true, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Header block: "Address: ..." line - drop the tab and two of the spaces
#    that sit right before "Website" (tab+5sp -> 3sp only, no tab).
# ---------------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.ClearFormatting()
$old2 = [char]9 + "     Website"
$new2 = "   Website"
$find2.Execute($old2, $false, $true, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. The single-column table's grid width shrinks along with the new right
#    margin (table is auto-width / centered, so its content width tracks the
#    usable page width).
# ---------------------------------------------------------------------------
$t = $d.Tables(1)
$t.Columns(1).Width = 431.5   # 8630 twips

# ---------------------------------------------------------------------------
# 4. Page margins: top 2127->1440 twips, right 1127->1800 twips (values are
#    expressed in points for PageSetup, 20 twips == 1 point).
# ---------------------------------------------------------------------------
$ps = $d.PageSetup
$ps.TopMargin = 72     # 1440 twips
$ps.RightMargin = 90   # 1800 twips

# ---------------------------------------------------------------------------
# 5. Bullet text: drop the stray leading space on "Use of D3, JavaScript...".
# ---------------------------------------------------------------------------
$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Execute(" Use of D3, JavaScript", $false, $true, $false, $false, $false, $true, 1, $false, "Use of D3, JavaScript", 2) | Out-Null

Write-Output "edits applied"
